$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert 8 new rows right before the existing row 36. This pushes the old
# rows 36-59 down to 44-67, matching the target layout.
# ---------------------------------------------------------------------------
$ws.Rows("36:43").Insert() | Out-Null

# ---------------------------------------------------------------------------
# Populate the 8 new rows (36-43) with the new Adsriver ad links.
# Column B gets the usual "ok" marker, column G gets the URL text + link.
# The list below is ordered the way the links were actually added to the
# workbook (controls the r:id numbering of the new relationships).
# ---------------------------------------------------------------------------
$newLinks = @(
    @{ Row = 36; Url = "http://www.adsriver.com/34/posts/20-Travel-Tickets/376--Vacation-Rentals/988030-Davao-Country-Side-Tour-.html" },
    @{ Row = 37; Url = "http://www.adsriver.com/34/posts/20-Travel-Tickets/377-Other-Travel-Ads/988031-Surigao-Tour-Package-.html" },
    @{ Row = 40; Url = "http://www.adsriver.com/34/posts/20-Travel-Tickets/377-Other-Travel-Ads/988032-Camiguin-Tour-Package.html" },
    @{ Row = 38; Url = "http://www.adsriver.com/34/posts/20-Travel-Tickets/377-Other-Travel-Ads/988034-Mati-Davao-Oriental-.html" },
    @{ Row = 39; Url = "http://www.adsriver.com/34/posts/20-Travel-Tickets/377-Other-Travel-Ads/988036-Davao-City-Tour.html" },
    @{ Row = 41; Url = "http://www.adsriver.com/34/posts/20-Travel-Tickets/377-Other-Travel-Ads/988037-Davao-White-Water-Rafting.html" },
    @{ Row = 42; Url = "http://www.adsriver.com/34/posts/20-Travel-Tickets/376--Vacation-Rentals/988040-Pearl-Farm-Beach-Resort.html" },
    @{ Row = 43; Url = "http://www.adsriver.com/34/posts/20-Travel-Tickets/376--Vacation-Rentals/988041-Samal-Island-Hopping-Package-.html" }
)

for ($r = 36; $r -le 43; $r++) {
    $ws.Range("B$r").Value = "ok"
}

foreach ($item in $newLinks) {
    $cell = $ws.Range("G$($item.Row)")
    $cell.Value = $item.Url
    $ws.Hyperlinks.Add($cell, $item.Url) | Out-Null
    # Re-apply the existing "Hyperlink" cell style (copied from another
    # hyperlinked cell) so we don't introduce a brand-new style record.
    $cell.Style = $ws.Range("G14").Style
}

# ---------------------------------------------------------------------------
# Update the view state to match (scroll position / active selection).
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 25
$ws.Range("C42").Select() | Out-Null
